$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add LoginDTO(id, password) / MemberDTO request-response entries
$ws.Range("G2").Value = "LoginDTO(id, password)"
$ws.Range("H2").Value = "MemberDTO"

# Row 4: change Method from GET to POST
$ws.Range("A4").Value = "POST"

# Row 11: add NoSuchElementException
$ws.Range("H11").Value = "NoSuchElementException"

# Row 4: add HttpServletResponse_OK response
$ws.Range("H4").Value = "HttpServletResponse_OK"

# Move the active selection to H5 as in the saved workbook
$ws.Range("H5").Select()
